$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-136 down to 55-137
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with its data
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = 'Vega Modelo de Temuco'
$ws.Range("C54").Value = 'La Araucanía'
$ws.Range("D54").Value = 44477
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = 100112005
$ws.Range("G54").Value = 'Puerro'
$ws.Range("H54").Value = 'Azul de Maquehue'
$ws.Range("I54").Value = 'Primera'
$ws.Range("J54").Value = 50
$ws.Range("K54").Value = 6000
$ws.Range("L54").Value = 7000
$ws.Range("M54").Value = 6600
$ws.Range("N54").Value = '$/docena de paquetes'
$ws.Range("O54").Value = 'Provincia de Cautín'
$ws.Range("P54").Value = 550
$ws.Range("Q54").Value = 12
$ws.Range("R54").Value = 'Hortaliza'
